{"js": "// The document contains five \"<id>p073r_N</id>\" tags (N = 1..5), each\n// originally split across three runs:\n//   run1 \"<id>\"       (Courier New, color 7f6000, sz 18)\n//   run2 \"p073r_N\"    (color 000000, default font)\n//   run3 \"</id>\"      (Courier New, color 7f6000, sz 18)\n// The edit merges each trio into a single run whose text is the full\n// \"<id>p073r_N</id>\" string, keeping the formatting of the first\n// (\"<id>\") run - this is what Word does when you select across the three\n// runs and retype/replace the text in one go.\nconst body = context.document.body;\n\nfor (let i = 1; i <= 5; i++) {\n  const idText = `<id>p073r_${i}</id>`;\n\n  // Locate the (currently three-run) occurrence of the full tag text -\n  // Word's search matches text across run boundaries.\n  const found = body.search(idText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  if (found.items.length > 0) {\n    // Replacing the found range's text with itself collapses the\n    // underlying runs into a single run that carries the first run's\n    // character formatting (Courier New / 7f6000 / 18pt).\n    found.items[0].insertText(idText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document contains five \"<id>p073r_N</id>\" tags (N = 1..5), each\n# originally split across three runs:\n#   run1 \"<id>\"       (Courier New, color 7f6000, sz 18)\n#   run2 \"p073r_N\"    (color 000000, default font)\n#   run3 \"</id>\"      (Courier New, color 7f6000, sz 18)\n#\n# The edit merges each trio back into a single run whose text is the\n# full \"<id>p073r_N</id>\" string, carrying the formatting of the first\n# (\"<id>\") run. A Find/Replace across the whole tag text (which spans\n# all three runs) collapses them into one run.\n$d = $word.ActiveDocument\n\nfor ($i = 1; $i -le 5; $i++) {\n    $idText = \"<id>p073r_$i</id>\"\n\n    $r = $d.Content\n    $r.Find.Execute(\n        $idText,   # FindText\n        $true,     # MatchCase\n        $false,    # MatchWholeWord\n        $false,    # MatchWildcards\n        $false,    # MatchSoundsLike\n        $false,    # MatchAllWordForms\n        $true,     # Forward\n        1,         # Wrap (wdFindContinue)\n        $false,    # Format\n        $idText,   # ReplaceWith\n        2          # Replace (wdReplaceAll)\n    ) | Out-Null\n}\n"}
